# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect the newly scraped counts.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAllTypes   = $wb.Worksheets.Item("全部类型")

# Row -> New value, for the "展览" sheet
$exhibitionUpdates = @{
    4  = 12205
    5  = 4490
    6  = 45
    9  = 26
    10 = 2606
    11 = 1130
    16 = 208
    18 = 11462
    19 = 11538
}

foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Row -> New value, for the "全部类型" sheet
$allTypesUpdates = @{
    4  = 12205
    5  = 4490
    6  = 45
    9  = 26
    10 = 2606
    12 = 1130
    17 = 208
    19 = 11462
    20 = 11538
}

foreach ($row in $allTypesUpdates.Keys) {
    $sheetAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
